$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Part 1: swap/rotate row data for re-ordered matches (columns F:V) ----
function Swap-Rows($ws, $rowA, $rowB) {
    $va = $ws.Range("F$rowA" + ":V$rowA").Value2()
    $vb = $ws.Range("F$rowB" + ":V$rowB").Value2()
    $ws.Range("F$rowA" + ":V$rowA").Value2 = $vb
    $ws.Range("F$rowB" + ":V$rowB").Value2 = $va
}

function Rotate-Rows($ws, [int[]]$rows) {
    # new content of rows[i] = old content of rows[i+1] (wrapping)
    $n = $rows.Count
    $orig = New-Object 'object[]' $n
    for ($i = 0; $i -lt $n; $i++) {
        $r = $rows[$i]
        $orig[$i] = $ws.Range("F$r" + ":V$r").Value2()
    }
    for ($i = 0; $i -lt $n; $i++) {
        $r = $rows[$i]
        $src = $orig[($i + 1) % $n]
        $ws.Range("F$r" + ":V$r").Value2 = $src
    }
}

Swap-Rows $ws 28 29
Swap-Rows $ws 47 48
Swap-Rows $ws 61 62
Swap-Rows $ws 63 64
Rotate-Rows $ws @(69, 71, 72)
Swap-Rows $ws 77 78
Swap-Rows $ws 93 94
Swap-Rows $ws 100 101
Rotate-Rows $ws @(109, 110, 111, 112)
Rotate-Rows $ws @(116, 119, 117, 118)

# ---- Part 2: append 7 new rows (124-130) with data + formatting ----
$ws.Range("A123:V123").Copy($ws.Range("A124:V124"))
$ws.Range("A123:V123").Copy($ws.Range("A125:V125"))
$ws.Range("A123:V123").Copy($ws.Range("A126:V126"))
$ws.Range("A123:V123").Copy($ws.Range("A127:V127"))
$ws.Range("A123:V123").Copy($ws.Range("A128:V128"))
$ws.Range("A123:V123").Copy($ws.Range("A129:V129"))
$ws.Range("A123:V123").Copy($ws.Range("A130:V130"))

$row124 = New-Object 'object[,]' 1,22
$row124[0,0] = 123.0
$row124[0,1] = "poland"
$row124[0,2] = "iii-liga-group-iii"
$row124[0,3] = "2023-2024"
$row124[0,4] = 45248.45833333334
$row124[0,5] = "Rakow II"
$row124[0,6] = 2.0
$row124[0,7] = "Zielona Gora"
$row124[0,8] = 2.0
$row124[0,9] = 2.86
$row124[0,10] = "18/11/2023 00:12"
$row124[0,11] = 3.46
$row124[0,12] = "18/11/2023 10:47"
$row124[0,13] = 3.51
$row124[0,14] = "18/11/2023 00:12"
$row124[0,15] = 3.09
$row124[0,16] = "18/11/2023 10:47"
$row124[0,17] = 2.1
$row124[0,18] = "18/11/2023 00:12"
$row124[0,19] = 2.03
$row124[0,20] = "18/11/2023 10:44"
$row124[0,21] = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/rks-rakow-czestochowa-zielona-gora/OIVZ8oVa/"
$ws.Range("A124:V124").Value2 = $row124

$row125 = New-Object 'object[,]' 1,22
$row125[0,0] = 124.0
$row125[0,1] = "poland"
$row125[0,2] = "iii-liga-group-iii"
$row125[0,3] = "2023-2024"
$row125[0,4] = 45248.5
$row125[0,5] = "Bytom Odrzanski"
$row125[0,6] = 1.0
$row125[0,7] = "Starowice Dolne"
$row125[0,8] = 0.0
$row125[0,9] = 1.61
$row125[0,10] = "18/11/2023 01:13"
$row125[0,11] = 1.69
$row125[0,12] = "18/11/2023 11:12"
$row125[0,13] = 4.0
$row125[0,14] = "18/11/2023 01:13"
$row125[0,15] = 3.74
$row125[0,16] = "18/11/2023 11:12"
$row125[0,17] = 3.96
$row125[0,18] = "18/11/2023 01:13"
$row125[0,19] = 4.05
$row125[0,20] = "18/11/2023 11:04"
$row125[0,21] = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/bytom-odrzanski-starowice-dolne/xM94g6VC/"
$ws.Range("A125:V125").Value2 = $row125

$row126 = New-Object 'object[,]' 1,22
$row126[0,0] = 125.0
$row126[0,1] = "poland"
$row126[0,2] = "iii-liga-group-iii"
$row126[0,3] = "2023-2024"
$row126[0,4] = 45248.54166666666
$row126[0,5] = "Gwarek Tarnowskie Gory"
$row126[0,6] = 3.0
$row126[0,7] = "Carina Gubin"
$row126[0,8] = 1.0
$row126[0,9] = 2.01
$row126[0,10] = "18/11/2023 02:12"
$row126[0,11] = 2.07
$row126[0,12] = "18/11/2023 12:54"
$row126[0,13] = 3.55
$row126[0,14] = "18/11/2023 02:12"
$row126[0,15] = 3.48
$row126[0,16] = "18/11/2023 12:54"
$row126[0,17] = 2.92
$row126[0,18] = "18/11/2023 02:12"
$row126[0,19] = 2.99
$row126[0,20] = "18/11/2023 12:52"
$row126[0,21] = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gwarek-tarnowskie-gory-carina-gubin/WKjeGTxQ/"
$ws.Range("A126:V126").Value2 = $row126

$row127 = New-Object 'object[,]' 1,22
$row127[0,0] = 126.0
$row127[0,1] = "poland"
$row127[0,2] = "iii-liga-group-iii"
$row127[0,3] = "2023-2024"
$row127[0,4] = 45248.54166666666
$row127[0,5] = "Polkowice"
$row127[0,6] = 1.0
$row127[0,7] = "Bielsko-Biala"
$row127[0,8] = 1.0
$row127[0,9] = 2.64
$row127[0,10] = "18/11/2023 01:13"
$row127[0,11] = 2.76
$row127[0,12] = "18/11/2023 12:54"
$row127[0,13] = 3.48
$row127[0,14] = "18/11/2023 01:13"
$row127[0,15] = 3.28
$row127[0,16] = "18/11/2023 12:54"
$row127[0,17] = 2.2
$row127[0,18] = "18/11/2023 01:13"
$row127[0,19] = 2.29
$row127[0,20] = "18/11/2023 12:54"
$row127[0,21] = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/polkowice-rekord-bielsko-biala/l0wNBmos/"
$ws.Range("A127:V127").Value2 = $row127

$row128 = New-Object 'object[,]' 1,22
$row128[0,0] = 127.0
$row128[0,1] = "poland"
$row128[0,2] = "iii-liga-group-iii"
$row128[0,3] = "2023-2024"
$row128[0,4] = 45248.5625
$row128[0,5] = "Gornik Zabrze II"
$row128[0,6] = 2.0
$row128[0,7] = "Sleza Wroclaw"
$row128[0,8] = 3.0
$row128[0,9] = 2.49
$row128[0,10] = "18/11/2023 02:42"
$row128[0,11] = 2.28
$row128[0,12] = "18/11/2023 08:15"
$row128[0,13] = 3.62
$row128[0,14] = "18/11/2023 02:42"
$row128[0,15] = 3.52
$row128[0,16] = "18/11/2023 11:34"
$row128[0,17] = 2.26
$row128[0,18] = "18/11/2023 02:42"
$row128[0,19] = 2.58
$row128[0,20] = "18/11/2023 08:15"
$row128[0,21] = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/gornik-zabrze-sleza-wroclaw/hUiiH9NJ/"
$ws.Range("A128:V128").Value2 = $row128

$row129 = New-Object 'object[,]' 1,22
$row129[0,0] = 128.0
$row129[0,1] = "poland"
$row129[0,2] = "iii-liga-group-iii"
$row129[0,3] = "2023-2024"
$row129[0,4] = 45248.5625
$row129[0,5] = "Pawlowice"
$row129[0,6] = 1.0
$row129[0,7] = "Stilon Gorzow"
$row129[0,8] = 0.0
$row129[0,9] = 1.83
$row129[0,10] = "18/11/2023 02:42"
$row129[0,11] = 1.86
$row129[0,12] = "18/11/2023 03:21"
$row129[0,13] = 3.78
$row129[0,14] = "18/11/2023 02:42"
$row129[0,15] = 3.72
$row129[0,16] = "18/11/2023 11:34"
$row129[0,17] = 3.22
$row129[0,18] = "18/11/2023 02:42"
$row129[0,19] = 3.28
$row129[0,20] = "18/11/2023 09:05"
$row129[0,21] = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/pniowek-pawlowice-stilon-gorzow/pSWV9RFg/"
$ws.Range("A129:V129").Value2 = $row129

$row130 = New-Object 'object[,]' 1,22
$row130[0,0] = 129.0
$row130[0,1] = "poland"
$row130[0,2] = "iii-liga-group-iii"
$row130[0,3] = "2023-2024"
$row130[0,4] = 45248.57291666666
$row130[0,5] = "Slask Wroclaw II"
$row130[0,6] = 3.0
$row130[0,7] = "Jelenia Gora"
$row130[0,8] = 0.0
$row130[0,9] = 1.6
$row130[0,10] = "18/11/2023 03:13"
$row130[0,11] = 1.68
$row130[0,12] = "18/11/2023 13:35"
$row130[0,13] = 4.07
$row130[0,14] = "18/11/2023 03:13"
$row130[0,15] = 4.06
$row130[0,16] = "18/11/2023 13:35"
$row130[0,17] = 3.98
$row130[0,18] = "18/11/2023 03:13"
$row130[0,19] = 3.81
$row130[0,20] = "18/11/2023 13:35"
$row130[0,21] = "https://www.betexplorer.com/football/poland/iii-liga-group-iii/slask-wroclaw-karkonosze-jelenia-gora/23sRA70m/"
$ws.Range("A130:V130").Value2 = $row130

"Dimension: " + $ws.UsedRange.Address()